# pre-final poster update (fixed enclosure section)
#
# Applies two changes to slide 1 of the WATech Park poster:
#   1. Nudge the "METHOD" section header rectangle (Rectangle 7) down slightly.
#   2. Reposition/resize/retitle the servo-horn-extender caption textbox
#      (TextBox 56) in the enclosure section.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Rectangle 7 ("METHOD" header box): a:off y 5181600 -> 5257800 EMU
#    (408pt -> 414pt). Left/width/height are untouched.
# ---------------------------------------------------------------------
$rect7 = $s.Shapes.Item("Rectangle 7")
$rect7.Top = 414.0

# ---------------------------------------------------------------------
# 2) TextBox 56 (servo barrier caption): move/widen the box and update
#    its caption text.
#      a:off  x 33985200 -> 33093315 EMU (2676pt      -> 2605.7728346456693pt)
#      a:off  y 16816630 -> 16833919 EMU (1324.144...pt -> 1325.505433070866pt)
#      a:ext cx 5300727 -> 5986527 EMU (417.3800787...pt -> 471.3800964355469pt)
#      a:ext cy 430887  -> 430887  EMU (unchanged)
# ---------------------------------------------------------------------
$caption = $s.Shapes.Item("TextBox 56")
$caption.Left = 2605.7728346456693
$caption.Top = 1325.505433070866
$caption.Width = 471.3800964355469
$caption.TextFrame.TextRange.Text = "3D Printed SG90 Servo Horn Extender(barrier)"
